$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mean_Intake (col O) and sem_Intake (col P): previous values were raw
# milligrams intaken; divide by animal weight to get the correct intake metric.
$ws.Range("O2").Value = 1013.3541529367591
$ws.Range("P2").Value = 223.50905325444523
$ws.Range("O3").Value = 1012.6725564179686
$ws.Range("P3").Value = 218.1014657261255
$ws.Range("O4").Value = 869.68670755655762
$ws.Range("P4").Value = 194.19713312888376
$ws.Range("O5").Value = 969.58536317097639
$ws.Range("P5").Value = 189.87261388109241
$ws.Range("O6").Value = 935.45722086368585
$ws.Range("P6").Value = 125.60844660689307
$ws.Range("O7").Value = 1058.4744834011026
$ws.Range("P7").Value = 148.9176572687704
$ws.Range("O8").Value = 1062.9633821786867
$ws.Range("P8").Value = 150.37190999404399
$ws.Range("O9").Value = 1307.6147416838887
$ws.Range("P9").Value = 235.0643201677463
$ws.Range("O10").Value = 1234.1926776994983
$ws.Range("P10").Value = 282.30160267793923
$ws.Range("O11").Value = 1214.0962080538557
$ws.Range("P11").Value = 282.77463501870756
$ws.Range("O12").Value = 1198.6011843622302
$ws.Range("P12").Value = 242.30707532578847
$ws.Range("O13").Value = 1193.7060900255703
$ws.Range("P13").Value = 212.97821622205674
$ws.Range("O14").Value = 1148.51923498107
$ws.Range("P14").Value = 186.24049534958002
$ws.Range("O15").Value = 1051.6128144942913
$ws.Range("P15").Value = 189.91703814034412
$ws.Range("O16").Value = 1058.6948163905847
$ws.Range("P16").Value = 157.32769425745403
$ws.Range("O17").Value = 135.8047085719584
$ws.Range("P17").Value = 47.440191316611013
$ws.Range("O18").Value = 82.110134038026217
$ws.Range("P18").Value = 20.276116855441522
$ws.Range("O19").Value = 242.09964734485072
$ws.Range("P19").Value = 104.88975574093729
$ws.Range("O20").Value = 560.25344571454104
$ws.Range("P20").Value = 166.94765708029666
$ws.Range("O21").Value = 629.66067425201345
$ws.Range("P21").Value = 193.85347490393204
$ws.Range("O22").Value = 1118.0876587960811
$ws.Range("P22").Value = 235.32063753393194
$ws.Range("O23").Value = 1194.6399077796943
$ws.Range("P23").Value = 226.82999973487986
$ws.Range("O24").Value = 1414.0979456422438
$ws.Range("P24").Value = 243.35371403633505
$ws.Range("O25").Value = 1874.2925638964371
$ws.Range("P25").Value = 463.43867495472182
$ws.Range("O26").Value = 1958.5594495516073
$ws.Range("P26").Value = 357.54820917167223
$ws.Range("O27").Value = 2011.6276273441558
$ws.Range("P27").Value = 297.21528318964096
$ws.Range("O28").Value = 1870.5022681164751
$ws.Range("P28").Value = 326.61111683185891
$ws.Range("O29").Value = 1620.7164167702267
$ws.Range("P29").Value = 276.73646169487245
$ws.Range("O30").Value = 1812.5157070507901
$ws.Range("P30").Value = 317.27077489965501
$ws.Range("O31").Value = 2244.9500767703648
$ws.Range("P31").Value = 453.30191787929903
$ws.Range("O32").Value = 2219.9478514996617
$ws.Range("P32").Value = 399.54866788205146
$ws.Range("O33").Value = 1014.889590936324
$ws.Range("P33").Value = 247.43505663383664
$ws.Range("O34").Value = 1023.5884829543371
$ws.Range("P34").Value = 402.48811133830623
$ws.Range("O35").Value = 1062.8143485577211
$ws.Range("P35").Value = 450.81110618322492
$ws.Range("O36").Value = 425.73456664521922
$ws.Range("P36").Value = 76.207963456188963
$ws.Range("O37").Value = 876.67788180511855
$ws.Range("P37").Value = 225.7027668121369
$ws.Range("O38").Value = 782.97013267340162
$ws.Range("P38").Value = 134.84781477564408
$ws.Range("O39").Value = 831.73760691404345
$ws.Range("P39").Value = 139.9392136028591
$ws.Range("O40").Value = 958.43905837389161
$ws.Range("P40").Value = 178.60495545824884
$ws.Range("O41").Value = 926.50754983630713
$ws.Range("P41").Value = 150.82524869402749
$ws.Range("O42").Value = 934.66152921890341
$ws.Range("P42").Value = 166.92776307696349
$ws.Range("O43").Value = 1107.30351296396
$ws.Range("P43").Value = 197.27741620435074
$ws.Range("O44").Value = 762.58153466419026
$ws.Range("P44").Value = 117.60484350308596
$ws.Range("O45").Value = 980.71342063600048
$ws.Range("P45").Value = 189.48058483926411
$ws.Range("O46").Value = 800.93681712347393
$ws.Range("P46").Value = 151.34599071903207
$ws.Range("O47").Value = 980.33997493778088
$ws.Range("P47").Value = 194.49117275704072
$ws.Range("O48").Value = 101.24748633857224
$ws.Range("P48").Value = 57.870559127304546
$ws.Range("O49").Value = 89.676990153019759
$ws.Range("P49").Value = 32.108626977563823
$ws.Range("O50").Value = 214.06929696928583
$ws.Range("P50").Value = 66.006443310897254
$ws.Range("O51").Value = 357.17993232893696
$ws.Range("P51").Value = 90.288751702313206
$ws.Range("O52").Value = 582.61607888282322
$ws.Range("P52").Value = 135.98669874227724
$ws.Range("O53").Value = 943.9061060348572
$ws.Range("P53").Value = 172.50070918918752
$ws.Range("O54").Value = 1138.6233377147028
$ws.Range("P54").Value = 181.00231186431409
$ws.Range("O55").Value = 1121.8494022703881
$ws.Range("P55").Value = 163.2014559472521
$ws.Range("O56").Value = 1207.2343277447801
$ws.Range("P56").Value = 161.24101688280814
$ws.Range("O57").Value = 1246.333417808037
$ws.Range("P57").Value = 180.50453544392664
$ws.Range("O58").Value = 1249.6318604273133
$ws.Range("P58").Value = 181.24951987204344
$ws.Range("O59").Value = 1253.5483593684608
$ws.Range("P59").Value = 175.57153507993027
$ws.Range("O60").Value = 1196.6943321937997
$ws.Range("P60").Value = 169.05559930872815
$ws.Range("O61").Value = 1441.5717055705918
$ws.Range("P61").Value = 236.54510942477313
$ws.Range("O62").Value = 1677.556916506891
$ws.Range("P62").Value = 313.66186037042803
$ws.Range("O63").Value = 1434.5951666099381
$ws.Range("P63").Value = 191.25694012240439

# Column P (16) got slightly narrower after the overlay/format cleanup.
$ws.Columns.Item(16).ColumnWidth = 10.833333333333334
